$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 894.3333
$ws.Range("J43").Value = 896
$ws.Range("L43").Value = 896
$ws.Range("N43").Value = -1034

$ws.Range("H86").Value = 1374.9286
$ws.Range("I86").Value = 1312.375
$ws.Range("J86").Value = 1458.3334
$ws.Range("K86").Value = 1312.375
$ws.Range("L86").Value = 1458.3334
$ws.Range("M86").Value = -189.375
$ws.Range("N86").Value = -3704.3334

$ws.Range("H87").Value = 24723.334
$ws.Range("I87").Value = 23680
$ws.Range("J87").Value = 24818.182
$ws.Range("K87").Value = 23680
$ws.Range("L87").Value = 24818.182
$ws.Range("M87").Value = -22432
$ws.Range("N87").Value = -27314.182

$ws.Range("H89").Value = 1374.9286
$ws.Range("I89").Value = 1312.375
$ws.Range("J89").Value = 1458.3334
$ws.Range("K89").Value = 6561.875
$ws.Range("L89").Value = 7291.666999999999
$ws.Range("M89").Value = -945.875
$ws.Range("N89").Value = -18523.667

$ws.Range("H90").Value = 24723.334
$ws.Range("I90").Value = 23680
$ws.Range("J90").Value = 24818.182
$ws.Range("K90").Value = 71040
$ws.Range("L90").Value = 74454.546
$ws.Range("M90").Value = -64800
$ws.Range("N90").Value = -86934.546

$ws.Range("H116").Value = 17303886
$ws.Range("I116").Value = 27675418
$ws.Range("J116").Value = 18000
$ws.Range("K116").Value = 27675418
$ws.Range("L116").Value = 18000
$ws.Range("M116").Value = -27671976
$ws.Range("N116").Value = -24884

$ws.Range("H129").Value = 949.4666999999999
$ws.Range("I129").Value = 351.6
$ws.Range("J129").Value = 992.17145
$ws.Range("K129").Value = 1054.8
$ws.Range("L129").Value = 2976.51435
$ws.Range("M129").Value = 3945.2
$ws.Range("N129").Value = -12976.51435

$ws.Range("H137").Value = 1566.3077
$ws.Range("I137").Value = 1588.8695
$ws.Range("J137").Value = 1393.3334
$ws.Range("K137").Value = 4766.6085
$ws.Range("L137").Value = 4180.0002
$ws.Range("M137").Value = -2216.6085
$ws.Range("N137").Value = -9280.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5014.1113
$ws.Range("I32").Value = 3442.611
$ws.Range("J32").Value = 8157.1113
$ws.Range("K32").Value = 3442.611
$ws.Range("L32").Value = 8157.1113
$ws.Range("M32").Value = -3155.611
$ws.Range("N32").Value = -8731.1113

$ws.Range("H61").Value = 2481.9167
$ws.Range("I61").Value = 1651.3334
$ws.Range("J61").Value = 4973.6665
$ws.Range("K61").Value = 1651.3334
$ws.Range("L61").Value = 4973.6665
$ws.Range("M61").Value = -1439.3334
$ws.Range("N61").Value = -5397.6665

$ws.Range("H74").Value = 8393.177
$ws.Range("I74").Value = 1445.6
$ws.Range("J74").Value = 60500
$ws.Range("K74").Value = 1445.6
$ws.Range("L74").Value = 60500
$ws.Range("M74").Value = -571.5999999999999
$ws.Range("N74").Value = -62248

$ws.Range("H77").Value = 8393.177
$ws.Range("I77").Value = 1445.6
$ws.Range("J77").Value = 60500
$ws.Range("K77").Value = 7228
$ws.Range("L77").Value = 302500
$ws.Range("M77").Value = -2860
$ws.Range("N77").Value = -311236

$ws.Range("H102").Value = 2764.875
$ws.Range("I102").Value = 2731.2856
$ws.Range("K102").Value = 2731.2856
$ws.Range("M102").Value = -1109.2856

$ws.Range("H122").Value = 16002
$ws.Range("I122").Value = 26078
$ws.Range("K122").Value = 78234
$ws.Range("M122").Value = -75784

$ws.Range("H132").Value = 3800.348
$ws.Range("I132").Value = 3380.4167
$ws.Range("K132").Value = 10141.2501
$ws.Range("M132").Value = -7611.250100000001

$ws.Range("H133").Value = 64932.668
$ws.Range("J133").Value = 64932.668
$ws.Range("L133").Value = 64932.668
$ws.Range("N133").Value = -69992.66800000001

$ws.Range("H136").Value = 2481.9167
$ws.Range("I136").Value = 1651.3334
$ws.Range("J136").Value = 4973.6665
$ws.Range("K136").Value = 4954.0002
$ws.Range("L136").Value = 14920.9995
$ws.Range("M136").Value = -2404.0002
$ws.Range("N136").Value = -20020.9995

$ws.Range("H139").Value = 48886
$ws.Range("J139").Value = 48886
$ws.Range("L139").Value = 48886
$ws.Range("N139").Value = -59166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 187.90909
$ws.Range("I22").Value = 187.90909
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 187.90909
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -14.90908999999999
$ws.Range("N22").ClearContents()

$ws.Range("H107").Value = 674
$ws.Range("I107").Value = 605.7143
$ws.Range("J107").Value = 833.3333
$ws.Range("K107").Value = 605.7143
$ws.Range("L107").Value = 833.3333
$ws.Range("M107").Value = 1314.2857
$ws.Range("N107").Value = -4673.3333

$ws.Range("H133").Value = 51000
$ws.Range("J133").Value = 51000
$ws.Range("L133").Value = 51000
$ws.Range("N133").Value = -61120

$ws.Range("H134").Value = 3206
$ws.Range("I134").Value = 2401.238
$ws.Range("J134").Value = 5620.2856
$ws.Range("K134").Value = 7203.714
$ws.Range("L134").Value = 16860.8568
$ws.Range("M134").Value = -4668.714
$ws.Range("N134").Value = -21930.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 816.6667
$ws.Range("I16").Value = 725
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 725
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -438
$ws.Range("N16").Value = -1574

$ws.Range("H31").Value = 4273.439
$ws.Range("I31").Value = 1869.5883
$ws.Range("K31").Value = 1869.5883
$ws.Range("M31").Value = -1574.5883

$ws.Range("H34").Value = 4273.439
$ws.Range("I34").Value = 1869.5883
$ws.Range("K34").Value = 1869.5883
$ws.Range("M34").Value = -1667.5883

$ws.Range("H58").Value = 1957.0834
$ws.Range("I58").Value = 1094.625
$ws.Range("J58").Value = 3682
$ws.Range("K58").Value = 1094.625
$ws.Range("L58").Value = 3682
$ws.Range("M58").Value = -891.625
$ws.Range("N58").Value = -4088

$ws.Range("H113").Value = 816.6667
$ws.Range("I113").Value = 725
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 725
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1445
$ws.Range("N113").Value = -5340

$ws.Range("H132").Value = 3606.4119
$ws.Range("I132").Value = 1894.25
$ws.Range("K132").Value = 5682.75
$ws.Range("M132").Value = -3152.75

$ws.Range("H134").Value = 2940.3513
$ws.Range("I134").Value = 2011.6786
$ws.Range("J134").Value = 5829.5557
$ws.Range("K134").Value = 6035.0358
$ws.Range("L134").Value = 17488.6671
$ws.Range("M134").Value = -3500.0358
$ws.Range("N134").Value = -22558.6671

$ws.Range("H136").Value = 1957.0834
$ws.Range("I136").Value = 1094.625
$ws.Range("J136").Value = 3682
$ws.Range("K136").Value = 3283.875
$ws.Range("L136").Value = 11046
$ws.Range("M136").Value = -733.875
$ws.Range("N136").Value = -16146

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2755.169
$ws.Range("I131").Value = 293.33334
$ws.Range("J131").Value = 2863.7793
$ws.Range("K131").Value = 880.0000200000001
$ws.Range("L131").Value = 8591.3379
$ws.Range("M131").Value = 4159.99998
$ws.Range("N131").Value = -18671.3379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7756.25
$ws.Range("I70").Value = 9090.909
$ws.Range("J70").Value = 4820
$ws.Range("K70").Value = 9090.909
$ws.Range("L70").Value = 4820
$ws.Range("M70").Value = -8820.909
$ws.Range("N70").Value = -5360

$ws.Range("H73").Value = 7756.25
$ws.Range("I73").Value = 9090.909
$ws.Range("J73").Value = 4820
$ws.Range("K73").Value = 9090.909
$ws.Range("L73").Value = 4820
$ws.Range("M73").Value = -8154.909
$ws.Range("N73").Value = -6692

$ws.Range("H80").Value = 2313.3044
$ws.Range("J80").Value = 2253
$ws.Range("L80").Value = 2253
$ws.Range("N80").Value = -4249

$ws.Range("H83").Value = 2313.3044
$ws.Range("I83").Value = 2319.0476
$ws.Range("J83").Value = 2253
$ws.Range("L83").Value = 11265
$ws.Range("N83").Value = -21249

$ws.Range("H113").Value = 2272
$ws.Range("I113").Value = 1590.125
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 1590.125
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = 579.875
$ws.Range("N113").Value = -9339.5

$ws.Range("H122").Value = 795836.6
$ws.Range("I122").Value = 856901.0600000001
$ws.Range("K122").Value = 2570703.18
$ws.Range("M122").Value = -2568253.18

$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws.Range("H138").Value = 62500
$ws.Range("J138").Value = 62500
$ws.Range("L138").Value = 62500
$ws.Range("N138").Value = -72780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2983.95
$ws.Range("I100").Value = 2993.3333
$ws.Range("J100").Value = 2982.2942
$ws.Range("K100").Value = 2993.3333
$ws.Range("L100").Value = 2982.2942
$ws.Range("M100").Value = -2452.3333
$ws.Range("N100").Value = -4064.2942

$ws.Range("H132").Value = 3728.125
$ws.Range("I132").Value = 2192.3125
$ws.Range("J132").Value = 6799.75
$ws.Range("K132").Value = 6576.9375
$ws.Range("L132").Value = 20399.25
$ws.Range("M132").Value = -4046.9375
$ws.Range("N132").Value = -25459.25

$ws.Range("H136").Value = 2817.2666
$ws.Range("I136").Value = 1874.9166
$ws.Range("J136").Value = 6586.6665
$ws.Range("K136").Value = 5624.7498
$ws.Range("L136").Value = 19759.9995
$ws.Range("M136").Value = -3074.7498
$ws.Range("N136").Value = -24859.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 35722220
$ws.Range("I132").Value = 55565136
$ws.Range("K132").Value = 166695408
$ws.Range("M132").Value = -166692878

$ws.Range("H136").Value = 15199132
$ws.Range("I136").Value = 23882050
$ws.Range("K136").Value = 71646150
$ws.Range("M136").Value = -71643600
